$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: Croatia - Gibraltar
$ws.Range("B2").Value = "Croatia ✓ - Gibraltar: 3:0"
$ws.Range("G2").Value = "✓"

# Row 3: Bukovyna Chernivtsi - Podillya Khmelnytskyi
$ws.Range("B3").Value = "Bukovyna Chernivtsi ✓ - Podillya Khmelnytskyi: 2:0"
$ws.Range("G3").Value = "✓"

# Row 4: Louisville City FC - Miami FC
$ws.Range("B4").Value = "Louisville City FC ✓ - Miami FC: 1:0"
$ws.Range("G4").Value = "✓"

# Row 5: Netherlands - Finland
$ws.Range("B5").Value = "Netherlands ✓ - Finland: 4:0"
$ws.Range("G5").Value = "✓"

# Row 6: Scotland - Belarus
$ws.Range("B6").Value = "Scotland ✓ - Belarus: 2:1"
$ws.Range("G6").Value = "✓"

# Row 7: Romania - Austria
$ws.Range("B7").Value = "Romania - Austria X: 1:0"
$ws.Range("G7").Value = "X"

# Row 8: Zambia - Niger
$ws.Range("B8").Value = "Zambia X - Niger: 0:1"
$ws.Range("G8").Value = "X"

# Row 9: Egypt - Guinea-Bissau
$ws.Range("B9").Value = "Egypt ✓ - Guinea-Bissau: 1:0"
$ws.Range("G9").Value = "✓"

# Row 10: CD Plaza Amador - CD Árabe Unido
$ws.Range("B10").Value = "CD Plaza Amador X - CD Árabe Unido: 0:2"
$ws.Range("G10").Value = "X"

# Row 11: Atlético Tembetary - Club Libertad Asunción
$ws.Range("B11").Value = "Atlético Tembetary - Club Libertad Asunción X: 2:1"
$ws.Range("G11").Value = "X"

# Row 12: Burkina Faso - Ethiopia
$ws.Range("B12").Value = "Burkina Faso ✓ - Ethiopia: 3:1"
$ws.Range("G12").Value = "✓"

# Row 13: Ghana - Comoros
$ws.Range("B13").Value = "Ghana ✓ - Comoros: 1:0"
$ws.Range("G13").Value = "✓"

# Row 14: Mali - Madagascar
$ws.Range("B14").Value = "Mali ✓ - Madagascar: 4:1"
$ws.Range("G14").Value = "✓"

# Row 15: Sanfrecce Hiroshima - Yokohama FC
$ws.Range("B15").Value = "Sanfrecce Hiroshima ✓ - Yokohama FC: 2:1"
$ws.Range("G15").Value = "✓"

# Row 16: PSS Sleman - Kendal Tornado FC
$ws.Range("B16").Value = "PSS Sleman ✓ - Kendal Tornado FC: 3:1"
$ws.Range("G16").Value = "✓"

# Row 17: Club Deportivo Guabirá - Club Aurora
$ws.Range("B17").Value = "Club Deportivo Guabirá ✓ - Club Aurora: 2:1"
$ws.Range("G17").Value = "✓"

# Row 18: Pars Jonoubi Jam - Fard Alborz (double space, no mark, no G cell)
$ws.Range("B18").Value = "Pars Jonoubi Jam  - Fard Alborz: 0:0"

# Row 19: Antigua GFC - Deportivo Achuapa (double space, no mark, no G cell)
$ws.Range("B19").Value = "Antigua GFC  - Deportivo Achuapa: 02:00"

# Row 20: CA Estudiantes - Club Deportivo Maipú
$ws.Range("B20").Value = "CA Estudiantes ✓ - Club Deportivo Maipú: 1:0"
$ws.Range("G20").Value = "✓"

# Row 21: Jeonnam Dragons - Ansan Greeners
$ws.Range("B21").Value = "Jeonnam Dragons X - Ansan Greeners: 0:1"
$ws.Range("G21").Value = "X"

# Row 22: CSD Municipal - Cobán Imperial (double space, no mark, no G cell)
$ws.Range("B22").Value = "CSD Municipal  - Cobán Imperial: 00:00"
